$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173257112503052
$ws.Range("B1").Value = 2.389649629592896
$ws.Range("D1").Value = 2.365973234176636
$ws.Range("E1").Value = 1.209512710571289
